# Poker League 2023 workbook update
# - Adds a new game week (column AH, played 2023-11-14) with results for the
#   players who took part: Ashish, Sid, Panos, Chris, Kartik, Tanish, Yufeng.
# - Re-sorts the league table by TABLES (column D, games played) descending,
#   which is the sheet's existing sort order, so Yufeng's row moves up as her
#   game count increases from 5 to 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New date header in AH1, matching the style of the other date cells ---
$ws.Range("AH1").NumberFormat = $ws.Range("AG1").NumberFormat
$ws.Range("AH1").Value = 45244

# --- 2. New week's results for the players who played ---
$ws.Range("AH2").Value = 27.25    # Ashish
$ws.Range("AH3").Value = 4.9      # Sid
$ws.Range("AH4").Value = -5.59    # Panos
$ws.Range("AH5").Value = -20      # Chris
$ws.Range("AH6").Value = -40      # Kartik
$ws.Range("AH7").Value = 23.2     # Tanish
$ws.Range("AH13").Value = 10.6    # Yufeng (still on row 13 prior to the re-sort)

# --- 3. Re-sort the table (A2:AH34) by TABLES (D4:D34) descending, same as
#        the existing sheet sort state, now extended to the new column ---
$sortRange = $ws.Range("A2:AH34")
$sortKey = $ws.Range("D4:D34")
$sortRange.Sort($sortKey, 2)

# Refresh the sheet's remembered sort settings so they cover the expanded
# range too.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("D4:D34"), 0, 2)
$sortObj.SetRange($ws.Range("A2:AH34"))
$sortObj.Header = -4142
$sortObj.Apply()

# --- 4. Restore the frozen-pane view / selection as left by the edit ---
$ws.Range("H6").Select()
